$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.561.09'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '3.605.29'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.82'
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.29'
$ws.Range("E6").Value = '  +3.32%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -0.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '8.09'
$ws.Range("E9").Value = '  +2.58%  '
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("D12").Value = '4.217.19'
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.85'
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").Value = '3.594.39'
$ws.Range("E15").Value = '  +1.18%  '
$ws.Range("D16").Value = '66.650.49'
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("E17").Value = '  +0.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.53'
$ws.Range("E18").Value = '  +1.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.40'
$ws.Range("E19").Value = '  +3.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.13'
$ws.Range("E20").Value = '  +2.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '428.33'
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.619'
$ws.Range("E22").Value = '  +1.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '78.86'
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("D24").Value = '3.754.46'
$ws.Range("E24").Value = '  +1.56%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  +4.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.34'
$ws.Range("E27").Value = '  +4.87%  '
$ws.Range("E28").Value = '  +4.02%  '
$ws.Range("E29").Value = '  +0.67%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").Value = '  +1.13%  '
$ws.Range("D32").Value = '3.602.99'
$ws.Range("E32").Value = '  +1.52%  '
$ws.Range("E33").Value = '  +4.06%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.88'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").Value = '  +1.37%  '
$ws.Range("E38").Value = '  -1.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '177.81'
$ws.Range("E39").Value = '  +1.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0858'
$ws.Range("E40").Value = '  +0.91%  '
$ws.Range("E41").Value = '  +0.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.901'
$ws.Range("E42").Value = '  +0.93%  '
$ws.Range("E43").Value = '  -0.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.58'
$ws.Range("E44").Value = '  +10.57%  '
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("E46").Value = '  -1.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.12'
$ws.Range("E47").Value = '  -2.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.09'
$ws.Range("E48").Value = '  +1.84%  '
$ws.Range("E49").Value = '  +1.34%  '
$ws.Range("E50").Value = '  +1.67%  '
$ws.Range("D51").Value = '2.432.11'
$ws.Range("E51").Value = '  +5.69%  '
